$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0
$ws.Range("E2").Value = 22.6600000000001
$ws.Range("G2").Value = 0.02748020294164433
$ws.Range("H2").Value = 0.04517928539833432
$ws.Range("I2").Value = "'"
$ws.Range("I2").Style = "Normal"
$ws.Range("K2").Value = 4.386803861447532
$ws.Range("L2").Value = "[0.02139438487797385, 8.752213338017091]"
$ws.Range("M2").Value = 0.04890275104014496
$ws.Range("N2").Value = 0.04890275104014496
$ws.Range("O2").Value = -1.169842309444848
$ws.Range("P2").Value = "[-2.528368862348542, 0.1886842434588467]"
$ws.Range("Q2").Value = 0.09089197372654456
$ws.Range("R2").Value = 0.1817839474530891
$ws.Range("S2").Value = 11.57591124250821
$ws.Range("T2").Value = "[9.305176132533914, 13.846646352482496]"
$ws.Range("W2").Value = 4.218978978978999
$ws.Range("X2").Value = -0.6804804804804823
$ws.Range("Y2").Value = 9.118438438438481

# Row 3
$ws.Range("B3").Value = 0
$ws.Range("E3").Value = 23.85000000000029
$ws.Range("G3").Value = [double]"5.41762125096934e-05"
$ws.Range("H3").Value = 0.0005980838859669067
$ws.Range("K3").Value = 5.192658876180018
$ws.Range("L3").Value = "[2.688103102208717, 7.697214650151319]"
$ws.Range("M3").Value = [double]"5.830105142079489e-05"
$ws.Range("N3").Value = 0.0001166021028415898
$ws.Range("O3").Value = -0.01257894956392303
$ws.Range("P3").Value = "[-0.6478159025420389, 0.6226580034141929]"
$ws.Range("Q3").Value = 0.9689362653859799
$ws.Range("R3").Value = 0.9689362653859799
$ws.Range("S3").Value = 11.24578940308216
$ws.Range("T3").Value = "[9.639403660102694, 12.852175146061626]"
$ws.Range("W3").Value = 0.04774774774774571
$ws.Range("X3").Value = -2.363513513513543
$ws.Range("Y3").Value = 2.459009009009034
